$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1965
$ws.Range("A1965").Value = "Allianz Football League Roinn 2"
$ws.Range("B1965").Value = "Cavan"
$ws.Range("C1965").Value = "Louth"
$ws.Range("D1965").Value = "28/02/2026"
$ws.Range("E1965").Value = "Kingspan Breffni, Cavan"
$ws.Range("F1965").Value = "1-12"
$ws.Range("G1965").Value = "1-22"

# Row 1966
$ws.Range("A1966").Value = "Allianz Football League Roinn 2"
$ws.Range("B1966").Value = "Tyrone"
$ws.Range("C1966").Value = "Offaly"
$ws.Range("D1966").Value = "28/02/2026"
$ws.Range("E1966").Value = "O Neill Park, Dungannon"
$ws.Range("F1966").Value = "0-28"
$ws.Range("G1966").Value = "2-13"

# Row 1967
$ws.Range("A1967").Value = "Allianz Football League Roinn 2"
$ws.Range("B1967").Value = "Kildare"
$ws.Range("C1967").Value = "Meath"
$ws.Range("D1967").Value = "28/02/2026"
$ws.Range("E1967").Value = "Cedral St Conleth's Newbridge"
$ws.Range("F1967").Value = "0-10"
$ws.Range("G1967").Value = "1-21"

# Row 1968
$ws.Range("A1968").Value = "Allianz Football League Roinn 3"
$ws.Range("B1968").Value = "Down"
$ws.Range("C1968").Value = "Fermanagh"
$ws.Range("D1968").Value = "28/02/2026"
$ws.Range("E1968").Value = "Páirc Esler, Newry"
$ws.Range("F1968").Value = "1-18"
$ws.Range("G1968").Value = "1-15"

# Row 1969
$ws.Range("A1969").Value = "Allianz Football League Roinn 3"
$ws.Range("B1969").Value = "Laois"
$ws.Range("C1969").Value = "Westmeath"
$ws.Range("D1969").Value = "28/02/2026"
$ws.Range("E1969").Value = "Laois Hire O'Moore Park"
$ws.Range("F1969").Value = "3-16"
$ws.Range("G1969").Value = "1-13"

# Row 1970
$ws.Range("A1970").Value = "Allianz Football League Roinn 3"
$ws.Range("B1970").Value = "Limerick"
$ws.Range("C1970").Value = "Clare"
$ws.Range("D1970").Value = "28/02/2026"
$ws.Range("E1970").Value = "Mick Neville Park Rathkeale"
$ws.Range("F1970").Value = "1-12"
$ws.Range("G1970").Value = "1-21"

# Row 1971
$ws.Range("A1971").Value = "Allianz Football League Roinn 4"
$ws.Range("B1971").Value = "Antrim"
$ws.Range("C1971").Value = "Wicklow"
$ws.Range("D1971").Value = "28/02/2026"
$ws.Range("E1971").Value = "Roger Casements, Portglenone"
$ws.Range("F1971").Value = "2-18"
$ws.Range("G1971").Value = "1-14"

# Row 1972
$ws.Range("A1972").Value = "Allianz Hurling League Roinn 1B"
$ws.Range("B1972").Value = "Carlow"
$ws.Range("C1972").Value = "Clare"
$ws.Range("D1972").Value = "28/02/2026"
$ws.Range("E1972").Value = "Netwatch Cullen Park"
$ws.Range("F1972").Value = "0-18"
$ws.Range("G1972").Value = "1-28"

# Row 1973
$ws.Range("A1973").Value = "Allianz Hurling League Roinn 2"
$ws.Range("B1973").Value = "Derry"
$ws.Range("C1973").Value = "London"
$ws.Range("D1973").Value = "28/02/2026"
$ws.Range("E1973").Value = "Find Insurance Celtic Park, Derry"
$ws.Range("F1973").Value = "0-14"
$ws.Range("G1973").Value = "3-7"

# Row 1974
$ws.Range("A1974").Value = "Allianz Hurling League Roinn 2"
$ws.Range("B1974").Value = "Mayo"
$ws.Range("C1974").Value = "Meath"
$ws.Range("D1974").Value = "28/02/2026"
$ws.Range("E1974").Value = "Tooreen (Adrian Freeman Park)"
$ws.Range("F1974").Value = "1-12"
$ws.Range("G1974").Value = "0-20"

# Row 1975
$ws.Range("A1975").Value = "Allianz Hurling League Roinn 3"
$ws.Range("B1975").Value = "Roscommon"
$ws.Range("C1975").Value = "Tyrone"
$ws.Range("D1975").Value = "28/02/2026"
$ws.Range("E1975").Value = "King & Moffatt Dr. Hyde Park"
$ws.Range("F1975").Value = "1-16"
$ws.Range("G1975").Value = "0-19"

# Row 1976
$ws.Range("A1976").Value = "Allianz Hurling League Roinn 4"
$ws.Range("B1976").Value = "Lancashire"
$ws.Range("C1976").Value = "Longford"
$ws.Range("D1976").Value = "28/02/2026"
$ws.Range("E1976").Value = "Abbottstown - GAA Centre of Excellence"
$ws.Range("F1976").Value = "0-12"
$ws.Range("G1976").Value = "3-18"

# Row 1977
$ws.Range("A1977").Value = "Allianz Hurling League Roinn 4"
$ws.Range("B1977").Value = "Leitrim"
$ws.Range("C1977").Value = "Cavan"
$ws.Range("D1977").Value = "28/02/2026"
$ws.Range("E1977").Value = "Heartland Credit Union Páirc Seán MacDiarmada"
$ws.Range("F1977").Value = "0-20"
$ws.Range("G1977").Value = "0-17"

# Row 1978
$ws.Range("A1978").Value = "Masita All-Ireland PPS Br Edmund Ignatius Rice Cup (Senior D Football)"
$ws.Range("B1978").Value = "Largy College Clones"
$ws.Range("C1978").Value = "Mountmellick Community School"
$ws.Range("D1978").Value = "28/02/2026"
# E1978 left blank (empty in source)
$ws.Range("F1978").Value = "5-15"
$ws.Range("G1978").Value = "2-1"

# Row 1979
$ws.Range("A1979").Value = "Masita All-Ireland PPS Br Edmund Ignatius Rice Cup (Senior D Football)"
$ws.Range("B1979").Value = "Rice College Ennis"
$ws.Range("C1979").Value = "St. Brendan's College, Belmullet"
$ws.Range("D1979").Value = "28/02/2026"
$ws.Range("E1979").Value = "Connacht GAA Centre"
$ws.Range("F1979").Value = "2-10"
$ws.Range("G1979").Value = "2-9"

# Row 1980
$ws.Range("A1980").Value = "Masita All-Ireland PPS Dr Eamonn O'Sullivan Cup (Senior C Football)"
$ws.Range("B1980").Value = "St Ciaran's, Ballygawley"
$ws.Range("C1980").Value = "St. Paul's College Raheny"
$ws.Range("D1980").Value = "28/02/2026"
$ws.Range("E1980").Value = "Crossmaglen"
$ws.Range("F1980").Value = "5-20"
$ws.Range("G1980").Value = "0-4"

# Row 1981
$ws.Range("A1981").Value = "Masita All-Ireland PPS Dr Eamonn O'Sullivan Cup (Senior C Football)"
$ws.Range("B1981").Value = "Balla Secondary School"
$ws.Range("C1981").Value = "Carrigaline Community School"
$ws.Range("D1981").Value = "28/02/2026"
$ws.Range("E1981").Value = "Clarecastle GAA"
$ws.Range("F1981").Value = "3-17"
$ws.Range("G1981").Value = "1-9"

# Row 1982
$ws.Range("A1982").Value = "Masita All-Ireland PPS Hogan Cup (Senior A Football)"
$ws.Range("B1982").Value = "St. Gerald's DLS College"
$ws.Range("C1982").Value = "Tralee CBS"
$ws.Range("D1982").Value = "28/02/2026"
$ws.Range("E1982").Value = "TUS Midwest"
$ws.Range("F1982").Value = "1-15"
$ws.Range("G1982").Value = "4-11"

# Row 1983
$ws.Range("A1983").Value = "Masita All-Ireland PPS Hogan Cup (Senior A Football)"
$ws.Range("B1983").Value = "Colaiste Mhuire Mullingar"
$ws.Range("C1983").Value = "Abbey Christian Brothers GS Newry"
$ws.Range("D1983").Value = "28/02/2026"
$ws.Range("E1983").Value = "Abbottstown - GAA Centre of Excellence"
$ws.Range("F1983").Value = "1-23"
$ws.Range("G1983").Value = "3-12"

# Row 1984
$ws.Range("A1984").Value = "Masita All-Ireland PPS Paddy Drummond Cup (Senior B Football)"
$ws.Range("B1984").Value = "Aquinas Diocesan GS Belfast"
$ws.Range("C1984").Value = "Cnoc Mhuire Granard"
$ws.Range("D1984").Value = "28/02/2026"
$ws.Range("E1984").Value = "Louis Leonard Memorial Park, Donagh"
$ws.Range("F1984").Value = "0-7"
$ws.Range("G1984").Value = "1-21"

# Row 1985
$ws.Range("A1985").Value = "Masita All-Ireland PPS Paddy Drummond Cup (Senior B Football)"
$ws.Range("B1985").Value = "Mount St. Michael Rosscarbery"
$ws.Range("C1985").Value = "St. Nathy's College Ballaghaderreen"
$ws.Range("D1985").Value = "28/02/2026"
$ws.Range("E1985").Value = "UL Grounds"
$ws.Range("F1985").Value = "0-11"
$ws.Range("G1985").Value = "3-18"
